$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph 1: opening sentence rewritten ---
Replace-Text `
    "Many years ago, there were three mighty gods who were planning the creation of a brand-new world that could keep their most unusual feats and creations. After that, they decided to add an unusual continent, so that all creations could live thus was created the continent of " `
    "Many years ago, three mighty gods planned the creation of a brand-new world that could keep their most unusual and creations. They decided to add an unusual continent, so that all creations could live, thus was created the continent of "

# --- Paragraph 2: rewritten entirely ---
Replace-Text `
    "The first inhabitants were the plans along with the animals, then the fishes, finally the humans and other creatures. All the gods were amazed at his creations. Therefore, they gathered and decided to gift their most valuable creations, the humans." `
    "The Gods created there most prized possession, humans. To help them thrive, plants and animals were added in abundance. As a final gift the Gods blessed them with their incredible powers."

# --- Paragraph 3: restructure god descriptions ---
# Delete the leading "The god of strength, " phrase (leaves just the tab in that run)
Replace-Text "The god of strength, " ""
Replace-Text ", has given strength to be dominant over adversity. The god of wisdom, " ", the God of strength, has given strength to be dominant over adversity. "
Replace-Text ", has bestowed pride for inner and mental growth along with intelligence so that they can master any teaching thus passed to them. Finally, the god of bravery, " ", the God of wisdom, has bestowed pride for inner and mental growth along with intelligence so that they can master any teaching passed unto them. Finally, "
Replace-Text ", granted the virtue that along with the vigor do not let themselves be so easily subdued or dominated." ", the God of bravery, granted the virtue that along with the vigor do not let themselves be so easily subdued or dominated."

# --- Paragraph 4: gifts / phrasing updates ---
Replace-Text `
    "After the bestowal of the divine presents, the gods ceased to exist, so that the humans could begin their journeys alone and discover the most unusual creations left by the vast continent of " `
    "After the bestowal of the divine gifts, the gods ceased to exist, so that humans could begin their journeys alone and discover the most unusual creations the Gods left on the vast and unknown continent of "
Replace-Text " still unknown." "."
